$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the M11 / N11 cell values (data edit from the commit)
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 1

# Move the view selection to O11 (matches the saved sheetView state)
$ws.Range("O11").Select()
